$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Reference/source cells (stable, unedited by this revision) used to
# --- carry over the correct style (number format) when a cell flips
# --- between a literal "N/A"-style string and a numeric value.
$srcStr0  = $ws.Cells.Item(14, 3)   # C14  s=14 t=s v=20 ("0")
$srcStrNA = $ws.Cells.Item(14, 5)   # E14  s=14 t=s v=21 ("***.*")
$srcNum15 = $ws.Cells.Item(30, 9)   # I30  s=15 (plain integer style)
$srcNum16 = $ws.Cells.Item(30, 11)  # K30  s=16 (signed 1-decimal style)

# --- Row 14 --------------------------------------------------------
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = -82.352941176470

# --- Row 15 --------------------------------------------------------
$srcStr0.Copy($ws.Cells.Item(15, 4))
$srcStrNA.Copy($ws.Cells.Item(15, 5))
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(15, 8).Value = -50
$ws.Cells.Item(15, 13).Value = -50

# --- Row 16 --------------------------------------------------------
$srcStr0.Copy($ws.Cells.Item(16, 4))
$srcStrNA.Copy($ws.Cells.Item(16, 5))
$ws.Cells.Item(16, 6).Value = 15
$ws.Cells.Item(16, 7).Value = 17
$ws.Cells.Item(16, 8).Value = -11.764705882352
$ws.Cells.Item(16, 9).Value = 81
$ws.Cells.Item(16, 11).Value = -2.409638554216
$ws.Cells.Item(16, 12).Value = 14.084507042253
$ws.Cells.Item(16, 13).Value = -19
$ws.Cells.Item(16, 14).Value = -79.230769230769

# --- Row 17 --------------------------------------------------------
$ws.Cells.Item(17, 4).Value = 8
$ws.Cells.Item(17, 5).Value = -62.5
$ws.Cells.Item(17, 6).Value = 16
$ws.Cells.Item(17, 7).Value = 15
$ws.Cells.Item(17, 8).Value = 6.666666666666
$ws.Cells.Item(17, 9).Value = 116
$ws.Cells.Item(17, 10).Value = 109
$ws.Cells.Item(17, 11).Value = 6.422018348623
$ws.Cells.Item(17, 12).Value = 1.754385964912
$ws.Cells.Item(17, 13).Value = 24.731182795698
$ws.Cells.Item(17, 14).Value = -65.060240963855

# --- Row 18 --------------------------------------------------------
$ws.Cells.Item(18, 6).Value = 11
$ws.Cells.Item(18, 7).Value = 8
$ws.Cells.Item(18, 8).Value = 37.5
$ws.Cells.Item(18, 9).Value = 60
$ws.Cells.Item(18, 10).Value = 85
$ws.Cells.Item(18, 11).Value = -29.411764705882
$ws.Cells.Item(18, 12).Value = -15.492957746478
$ws.Cells.Item(18, 13).Value = 5.263157894736
$ws.Cells.Item(18, 14).Value = -86.175115207373

# --- Row 19 --------------------------------------------------------
$ws.Cells.Item(19, 3).Value = 8
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 5).Value = -11.111111111111
$ws.Cells.Item(19, 6).Value = 27
$ws.Cells.Item(19, 7).Value = 31
$ws.Cells.Item(19, 8).Value = -12.903225806451
$ws.Cells.Item(19, 9).Value = 165
$ws.Cells.Item(19, 10).Value = 178
$ws.Cells.Item(19, 11).Value = -7.303370786516
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 13.793103448275
$ws.Cells.Item(19, 14).Value = -16.666666666666

# --- Row 20 --------------------------------------------------------
$ws.Cells.Item(20, 3).Value = 5
$ws.Cells.Item(20, 6).Value = 10
$ws.Cells.Item(20, 7).Value = 3
$ws.Cells.Item(20, 8).Value = 233.333333333333
$ws.Cells.Item(20, 9).Value = 32
$ws.Cells.Item(20, 11).Value = 39.130434782608
$ws.Cells.Item(20, 12).Value = 10.344827586206
$ws.Cells.Item(20, 13).Value = 255.555555555556
$ws.Cells.Item(20, 14).Value = -53.623188405797

# --- Row 21 --------------------------------------------------------
$ws.Cells.Item(21, 3).Value = 21
$ws.Cells.Item(21, 4).Value = 18
$ws.Cells.Item(21, 5).Value = 16.666666666666
$ws.Cells.Item(21, 6).Value = 80
$ws.Cells.Item(21, 7).Value = 76
$ws.Cells.Item(21, 8).Value = 5.263157894736
$ws.Cells.Item(21, 9).Value = 460
$ws.Cells.Item(21, 10).Value = 484
$ws.Cells.Item(21, 11).Value = -4.958677685950
$ws.Cells.Item(21, 12).Value = 1.098901098901
$ws.Cells.Item(21, 13).Value = 11.380145278450
$ws.Cells.Item(21, 14).Value = -68.406593406593

# --- Row 22 --------------------------------------------------------
$srcStr0.Copy($ws.Cells.Item(22, 4))
$srcStrNA.Copy($ws.Cells.Item(22, 5))
$srcStr0.Copy($ws.Cells.Item(22, 6))
$ws.Cells.Item(22, 7).Value = 2
$ws.Cells.Item(22, 8).Value = -100

# --- Row 23 --------------------------------------------------------
$ws.Cells.Item(23, 3).Value = 3
$srcNum15.Copy($ws.Cells.Item(23, 4))
$ws.Cells.Item(23, 4).Value = 3
$srcNum16.Copy($ws.Cells.Item(23, 5))
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 7
$ws.Cells.Item(23, 7).Value = 8
$ws.Cells.Item(23, 8).Value = -12.5
$ws.Cells.Item(23, 9).Value = 40
$ws.Cells.Item(23, 10).Value = 49
$ws.Cells.Item(23, 11).Value = -18.367346938775
$ws.Cells.Item(23, 12).Value = 2.564102564102
$ws.Cells.Item(23, 13).Value = 166.666666666667

# --- Row 24 --------------------------------------------------------
$ws.Cells.Item(24, 3).Value = 38
$ws.Cells.Item(24, 4).Value = 28
$ws.Cells.Item(24, 5).Value = 35.714285714285
$ws.Cells.Item(24, 6).Value = 110
$ws.Cells.Item(24, 7).Value = 104
$ws.Cells.Item(24, 8).Value = 5.769230769230
$ws.Cells.Item(24, 9).Value = 617
$ws.Cells.Item(24, 10).Value = 549
$ws.Cells.Item(24, 11).Value = 12.386156648451
$ws.Cells.Item(24, 12).Value = 13.419117647058
$ws.Cells.Item(24, 13).Value = 36.203090507726

# --- Row 25 --------------------------------------------------------
$ws.Cells.Item(25, 3).Value = 5
$ws.Cells.Item(25, 4).Value = 11
$ws.Cells.Item(25, 5).Value = -54.545454545454
$ws.Cells.Item(25, 6).Value = 31
$ws.Cells.Item(25, 7).Value = 36
$ws.Cells.Item(25, 8).Value = -13.888888888888
$ws.Cells.Item(25, 9).Value = 202
$ws.Cells.Item(25, 10).Value = 211
$ws.Cells.Item(25, 11).Value = -4.265402843601
$ws.Cells.Item(25, 12).Value = 2.020202020202
$ws.Cells.Item(25, 13).Value = -20.472440944881

# --- Row 26 --------------------------------------------------------
$srcStr0.Copy($ws.Cells.Item(26, 4))
$srcStrNA.Copy($ws.Cells.Item(26, 5))
$ws.Cells.Item(26, 7).Value = 3
$ws.Cells.Item(26, 8).Value = -33.333333333333

# --- Row 27 --------------------------------------------------------
$srcStr0.Copy($ws.Cells.Item(27, 3))
$srcNum15.Copy($ws.Cells.Item(27, 4))
$ws.Cells.Item(27, 4).Value = 1
$srcNum16.Copy($ws.Cells.Item(27, 5))
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 6).Value = 2
$ws.Cells.Item(27, 7).Value = 3
$ws.Cells.Item(27, 8).Value = -33.333333333333
$ws.Cells.Item(27, 9).Value = 19
$ws.Cells.Item(27, 10).Value = 28
$ws.Cells.Item(27, 11).Value = -32.142857142857
$ws.Cells.Item(27, 12).Value = 5.555555555555

# --- Row 28 --------------------------------------------------------
$srcNum15.Copy($ws.Cells.Item(28, 4))
$ws.Cells.Item(28, 4).Value = 2
$srcNum16.Copy($ws.Cells.Item(28, 5))
$ws.Cells.Item(28, 5).Value = -100
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(28, 10).Value = 9
$ws.Cells.Item(28, 11).Value = -33.333333333333
$ws.Cells.Item(28, 14).Value = -88

# --- Row 29 --------------------------------------------------------
$srcNum15.Copy($ws.Cells.Item(29, 4))
$ws.Cells.Item(29, 4).Value = 2
$srcNum16.Copy($ws.Cells.Item(29, 5))
$ws.Cells.Item(29, 5).Value = -100
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(29, 10).Value = 9
$ws.Cells.Item(29, 11).Value = -44.444444444444
$ws.Cells.Item(29, 14).Value = -88.888888888888

